$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 48, shifting the existing row 48 (and everything
# below it) down by one. This mirrors a new daily price record being added to
# the weekly consolidation for "Espárragos" at Vega Central Mapocho de Santiago.
$ws.Rows("48:48").Insert()

$ws.Range("A48").Value = 9
$ws.Range("B48").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = (Get-Date -Year 2022 -Month 4 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = 300000000
$ws.Range("G48").Value = "Espárragos"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Tercera"
$ws.Range("J48").Value = 7
$ws.Range("K48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("M48").Value = 30000
$ws.Range("N48").Value = "`$/bandeja 10 kilos"
$ws.Range("O48").Value = "Provincia de Linares"
$ws.Range("P48").Value = 3000
$ws.Range("Q48").Value = 10
$ws.Range("R48").Value = "Hortaliza"
